$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1): new "Wins" / "Losses" / "Ties" columns ---
$ws.Range("AD1").Value = "Wins"
$ws.Range("AE1").Value = "Losses"
$ws.Range("AF1").Value = "Ties"

# Match the existing bold/centered/bordered header style (same as A1:AC1)
# by copying the formatting from an existing header cell.
$ws.Range("A1").Copy()
$ws.Range("AD1:AF1").PasteSpecial(-4122)

# --- Data rows (2-51): New York Yankees' 2007 season record ---
# 94 wins, 68 losses, 0 ties - applied to every player row.
$ws.Range("AD2:AD51").Value = 94
$ws.Range("AE2:AE51").Value = 68
$ws.Range("AF2:AF51").Value = 0
